$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 7083.1665
$ws.Range("J43").Value = 8125
$ws.Range("L43").Value = 8125
$ws.Range("N43").Value = -8263

$ws.Range("H92").Value = 15625690
$ws.Range("I92").Value = 18519220
$ws.Range("K92").Value = 18519220
$ws.Range("M92").Value = -18517972

$ws.Range("H112").Value = 1496070.4
$ws.Range("J112").Value = 1638143.8
$ws.Range("L112").Value = 4914431.4
$ws.Range("N112").Value = -4916647.4

$ws.Range("H113").Value = 4420.8604
$ws.Range("I113").Value = 5365.8667
$ws.Range("K113").Value = 5365.8667
$ws.Range("M113").Value = -2111.8667

$ws.Range("H138").Value = 1780.78
$ws.Range("I138").Value = 820.53845
$ws.Range("J138").Value = 2394.7048
$ws.Range("K138").Value = 2461.61535
$ws.Range("L138").Value = 7184.1144
$ws.Range("M138").Value = 2678.38465
$ws.Range("N138").Value = -17464.1144

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H42").Value = 12842.667
$ws.Range("J42").Value = 14750
$ws.Range("L42").Value = 14750
$ws.Range("N42").Value = -15722

$ws.Range("H44").Value = 69993.664
$ws.Range("J44").Value = 69993.664
$ws.Range("L44").Value = 69993.664
$ws.Range("N44").Value = -70969.664

$ws.Range("H46").Value = 3617
$ws.Range("J46").Value = 3617
$ws.Range("L46").Value = 3617
$ws.Range("N46").Value = -4255

$ws.Range("H55").Value = 60007.75
$ws.Range("J55").Value = 69994.336
$ws.Range("L55").Value = 69994.336
$ws.Range("N55").Value = -70624.336

$ws.Range("H74").Value = 2325.7646
$ws.Range("I74").Value = 1927.9788
$ws.Range("J74").Value = 6999.75
$ws.Range("K74").Value = 1927.9788
$ws.Range("L74").Value = 6999.75
$ws.Range("M74").Value = -1053.9788
$ws.Range("N74").Value = -8747.75

$ws.Range("H77").Value = 2325.7646
$ws.Range("I77").Value = 1927.9788
$ws.Range("J77").Value = 6999.75
$ws.Range("K77").Value = 9639.894
$ws.Range("L77").Value = 34998.75
$ws.Range("M77").Value = -5271.894
$ws.Range("N77").Value = -43734.75

$ws.Range("H97").Value = 1302.8948
$ws.Range("I97").Value = 984.93335
$ws.Range("J97").Value = 2495.25
$ws.Range("K97").Value = 984.93335
$ws.Range("L97").Value = 2495.25
$ws.Range("M97").Value = -488.93335
$ws.Range("N97").Value = -3487.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3972173.2
$ws.Range("I134").Value = 5496478
$ws.Range("K134").Value = 16489434
$ws.Range("M134").Value = -16486899

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5188.6274
$ws.Range("I31").Value = 2542.8096
$ws.Range("J31").Value = 7040.7
$ws.Range("K31").Value = 2542.8096
$ws.Range("L31").Value = 7040.7
$ws.Range("M31").Value = -2247.8096
$ws.Range("N31").Value = -7630.7

$ws.Range("H34").Value = 5188.6274
$ws.Range("I34").Value = 2542.8096
$ws.Range("J34").Value = 7040.7
$ws.Range("K34").Value = 2542.8096
$ws.Range("L34").Value = 7040.7
$ws.Range("M34").Value = -2340.8096
$ws.Range("N34").Value = -7444.7

$ws.Range("H132").Value = 3332.4055
$ws.Range("I132").Value = 2924.3103
$ws.Range("J132").Value = 4811.75
$ws.Range("K132").Value = 8772.930899999999
$ws.Range("L132").Value = 14435.25
$ws.Range("M132").Value = -6242.930899999999
$ws.Range("N132").Value = -19495.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 113712250
$ws.Range("J4").Value = 199871860
$ws.Range("L4").Value = 599615580
$ws.Range("N4").Value = -599615804

$ws.Range("H9").Value = 4345.5
$ws.Range("I9").Value = 2500
$ws.Range("J9").Value = 4960.6665
$ws.Range("K9").Value = 7500
$ws.Range("L9").Value = 14881.9995
$ws.Range("M9").Value = -7276
$ws.Range("N9").Value = -15329.9995

$ws.Range("H34").Value = 455.42856
$ws.Range("J34").Value = 2665
$ws.Range("L34").Value = 7995
$ws.Range("N34").Value = -8163

$ws.Range("H39").Value = 4866.727
$ws.Range("J39").Value = 4866.727
$ws.Range("L39").Value = 14600.181
$ws.Range("N39").Value = -15188.181

$ws.Range("H55").Value = 3512.8462
$ws.Range("I55").Value = 1127.4
$ws.Range("J55").Value = 5003.75
$ws.Range("K55").Value = 3382.2
$ws.Range("L55").Value = 15011.25
$ws.Range("M55").Value = -3205.2
$ws.Range("N55").Value = -15365.25

$ws.Range("H131").Value = 1509.2858
$ws.Range("I131").Value = 626
$ws.Range("J131").Value = 2000
$ws.Range("K131").Value = 1878
$ws.Range("L131").Value = 6000
$ws.Range("M131").Value = 3162
$ws.Range("N131").Value = -16080

$ws.Range("H137").Value = 7657.95
$ws.Range("J137").Value = 9182.4375
$ws.Range("L137").Value = 27547.3125
$ws.Range("N137").Value = -37747.3125

$ws.Range("H139").Value = 3091.7778
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 14372.625
$ws.Range("J113").Value = 51205.5
$ws.Range("L113").Value = 51205.5
$ws.Range("N113").Value = -55545.5

$ws.Range("H120").Value = 79898.5
$ws.Range("J120").Value = 79898.5
$ws.Range("L120").Value = 79898.5
$ws.Range("N120").Value = -89574.5

$ws.Range("H122").Value = 1000
$ws.Range("I122").Value = 1000
$ws.Range("K122").Value = 3000
$ws.Range("M122").Value = -550

$ws.Range("H126").Value = 2739
$ws.Range("I126").Value = 2543.5715
$ws.Range("K126").Value = 7630.7145
$ws.Range("M126").Value = -5160.7145

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H112").Value = 96346.5
$ws.Range("J112").Value = 96346.5
$ws.Range("L112").Value = 96346.5
$ws.Range("N112").Value = -99300.5

$ws.Range("H122").Value = 14128.5
$ws.Range("I122").Value = 12685.728
$ws.Range("K122").Value = 38057.18399999999
$ws.Range("M122").Value = -35607.18399999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 18199.5
$ws.Range("I51").Value = 15399.667
$ws.Range("J51").Value = 20999.334
$ws.Range("K51").Value = 15399.667
$ws.Range("L51").Value = 20999.334
$ws.Range("M51").Value = -14889.667
$ws.Range("N51").Value = -22019.334

$ws.Range("H122").Value = 31253778
$ws.Range("I122").Value = 40004160
$ws.Range("J122").Value = 2410.1428
$ws.Range("K122").Value = 120012480
$ws.Range("L122").Value = 7230.428400000001
$ws.Range("M122").Value = -120010030
$ws.Range("N122").Value = -12130.4284
